{"js": "// The diff removes four paragraphs that sat between the\n// \"LOB1036: Geometria Anal\u00edtica (Requisito fraco)\" requirement line and\n// the trailing empty / page-break paragraphs at the very end of the\n// document body:\n//   1. an empty paragraph\n//   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3. another empty paragraph\n//   4. an empty paragraph with pageBreakBefore (+ left justification)\n//\n// Anchor on the unique \"Ver no Jupiter ...\" text and remove it together\n// with its immediate paragraph neighbours (one before, two after) so the\n// edit does not depend on hard-coded paragraph indices.\n\nconst body = context.document.body;\nconst results = body.search(\"Ver no Jupiter Salvar em pdf Salvar em docx\", {\n  matchCase: true,\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Anchor text \"Ver no Jupiter Salvar em pdf Salvar em docx\" not found.');\n}\n\nconst anchorParagraph = results.items[0].paragraphs.getFirst();\nconst beforeParagraph = anchorParagraph.getPrevious();\nconst afterParagraph1 = anchorParagraph.getNext();\nconst afterParagraph2 = afterParagraph1.getNext();\n\n// Delete the anchor paragraph and its three surrounding blank paragraphs.\nbeforeParagraph.delete();\nafterParagraph2.delete();\nafterParagraph1.delete();\nanchorParagraph.delete();\n\nawait context.sync();\n", "ps1": "# The diff removes four paragraphs that sat between the\n# \"LOB1036: Geometria Anal\u00edtica (Requisito fraco)\" requirement line and\n# the trailing empty / page-break paragraphs at the very end of the\n# document body:\n#   1. an empty paragraph\n#   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3. another empty paragraph\n#   4. an empty paragraph with pageBreakBefore (+ left justification)\n#\n# Locate the unique \"Ver no Jupiter ...\" paragraph and remove it together\n# with its immediate paragraph neighbours (one before, two after) so the\n# edit does not depend on hard-coded paragraph indices.\n\n$d = $word.ActiveDocument\n$needle = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n\n$anchor = $null\n$paragraphs = $d.Content.Paragraphs\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    $candidate = $paragraphs.Item($i)\n    if ($candidate.Range.Text.Contains($needle)) {\n        $anchor = $candidate\n        break\n    }\n}\n\nif ($anchor -eq $null) {\n    throw \"Anchor paragraph '$needle' not found\"\n}\n\n$before = $anchor.Previous()\n$after1 = $anchor.Next()\n$after2 = $after1.Next()\n\n# Delete from the last paragraph in document order back to the first so\n# earlier objects stay valid while later ones are removed.\n$after2.Range.Delete()\n$after1.Range.Delete()\n$anchor.Range.Delete()\n$before.Range.Delete()\n"}
